$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("For plotting")

$ws.Range("D2").Value = 0.0461968780928217
$ws.Range("E2").Value = 0.136094923708313

$ws.Range("D3").Value = 0.0453004886239722
$ws.Range("E3").Value = 0.148667067653778

$ws.Range("D4").Value = 0.0712061002463234
$ws.Range("E4").Value = 0.185263325055329
